$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.351.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.941.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7256'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3358'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.75'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07314'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8190'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08137'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.938.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.550'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.356.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008303'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '255.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.894'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.191.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9987'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.969'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.872'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.421'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1336'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.561'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.49%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.466'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.261'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05251'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.280'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7575'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.741'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02002'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.843'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.690'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '80.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4572'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.037'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8434'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.843'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.446'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.508'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4167'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.87%  '
